$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column E ("bg_auto"), shifting the existing
# bg_auto..lowess_span columns one place to the right. Excel carries the
# formatting of the column immediately to the left (D, "rel_quant") onto
# the freshly inserted column, which matches the header/value styling
# used by the new "quantifier_sel" column in the target workbook.
$ws.Range("E1").EntireColumn.Insert()

# New header and value for the inserted column.
$ws.Range("E1").Value = "quantifier_sel"
$ws.Range("E2").Value = 0

# Give the new column the same width as its left neighbor (D).
$ws.Range("E1").EntireColumn.ColumnWidth = $ws.Range("D1").EntireColumn.ColumnWidth

# Update the active selection to match the target workbook state.
$ws.Range("G12").Select()
